# Fix bad formatting: rows 27-42 on the "Categories" sheet had their
# data (columns B-G: key, description, externalId, name, slug, parent.key)
# shifted by a row. Re-sort the rows back into the correct order while
# keeping columns A (data-object) and H (parent.typeId) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Categories")

$data = @{
    27 = @("child-anyKey","child-anyDescription","child-anyId","Child-any","child-anySlug","ageGroupKey")
    28 = @("JonesKey","JonesDescription","JonesId","Jones","JonesSlug","brandKey")
    29 = @("BataleonKey","BataleonDescription","BataleonId","Bataleon","BataleonSlug","brandKey")
    30 = @("child-boyKey","child-boyDescription","child-boyId","Child-boy","child-boySlug","ageGroupKey")
    31 = @("SalomonKey","SalomonDescription","SalomonId","Salomon","SalomonSlug","brandKey")
    32 = @("RomeKey","RomeDescription","RomeId","Rome","RomeSlug","brandKey")
    33 = @("SimsKey","SimsDescription","SimsId","Sims","SimsSlug","brandKey")
    34 = @("SplitboardingKey","SplitboardingDescription","SplitboardingId","Splitboarding","SplitboardingSlug","terrainKey")
    35 = @("NitroKey","NitroDescription","NitroId","Nitro","NitroSlug","brandKey")
    36 = @("NeverSummerKey","NeverSummerDescription","NeverSummerId","NeverSummer","NeverSummerSlug","brandKey")
    37 = @("child-girlKey","child-girlDescription","child-girlId","Child-girl","child-girlSlug","ageGroupKey")
    38 = @("UnitedShapesKey","UnitedShapesDescription","UnitedShapesId","UnitedShapes","UnitedShapesSlug","brandKey")
    39 = @("3YearsKey","3YearsDescription","3YearsId","3Years","3YearsSlug","terrainKey")
    40 = @("YesKey","YesDescription","YesId","Yes","YesSlug","brandKey")
    41 = @("WestonKey","WestonDescription","WestonId","Weston","WestonSlug","brandKey")
    42 = @("CardiffKey","CardiffDescription","CardiffId","Cardiff","CardiffSlug","brandKey")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
